$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.956.21"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "1.893.37"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +1.60%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.40%  "
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4694"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3927"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.51"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08041"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.021"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "1.895.57"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.969"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.118"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06802"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.11%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001050"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.016"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("D22").Value = "27.977.97"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.516"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.346"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").Value = "2.123.24"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.461"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9721"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09515"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.667"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.401"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.360"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06135"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02257"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.219"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.110"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5994"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.270"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5703"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.403"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.939"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("B49").Value = "PaxosStandard"
$ws.Range("C49").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.112"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.56%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06934"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.62%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.99%  "
